# Applies the "add Jurisdiction row + refresh Date" edit to the Metadata sheet
# of the ValueSet workbook, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Refresh the "Date" metadata value (row 8, column B).
$ws.Cells.Item(8, 2).Value = "2024-07-01T07:50:29+00:00"

# 2. Insert a new "Jurisdiction" / "" row right before "Description" (currently row 11),
#    pushing Description / Purpose / Copyright / Immutable down by one row each
#    (11->12, 12->13, 13->14, 14->15), while preserving the exact original cell
#    styling (style index 2) for every moved / new cell.

# Remember the values currently in rows 11..14 (A and B columns) before we move them.
$vals = @()
for ($r = 11; $r -le 14; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $vals += ,@($a, $b)
}

# Row 15 doesn't exist yet, so copy the formatting of row 14 onto it first, which
# guarantees it reuses the existing "data row" style instead of Excel fabricating
# a brand-new (and subtly different) style entry.
$ws.Cells.Item(14, 1).Copy()
$ws.Range("A15").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(14, 2).Copy()
$ws.Range("B15").PasteSpecial(-4122)  # xlPasteFormats

# Shift the captured values down one row (12..15). Where the source cell had no
# value, clear the destination's contents (but keep its style).
for ($i = 0; $i -lt 4; $i++) {
    $destRow = 12 + $i
    $ws.Cells.Item($destRow, 1).Value = $vals[$i][0]
    if ($vals[$i][1]) {
        $ws.Cells.Item($destRow, 2).Value = $vals[$i][1]
    } else {
        $ws.Cells.Item($destRow, 2).ClearContents()
    }
}

# Finally, write the new "Jurisdiction" row into row 11 (it already carries the
# original "Description" row's style, which is style index 2).
$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""

Write-Host "Edit applied"
